$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TestData sheet previously had two "holes" in its 3-column table
# (A2 and B3 were never populated). The updated TestData fills them in
# with new values, reusing the same row style as the rest of the table.

$ws.Range("A2").Value = "abcd"
$ws.Range("B3").Value = "efgh"

# Make sure the newly-populated cells pick up the same cell formatting
# (style index) already used throughout the rest of the table.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
